$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: updated "Taxa" values
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 33
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 33
# B11..B15 unchanged (stay 0)

# Column C: updated "Data de Salvamento" timestamps for all data rows (2..15)
$newTimestamp = "2025-04-04 13:26:56"
$ws.Range("C2").Value = $newTimestamp
$ws.Range("C3").Value = $newTimestamp
$ws.Range("C4").Value = $newTimestamp
$ws.Range("C5").Value = $newTimestamp
$ws.Range("C6").Value = $newTimestamp
$ws.Range("C7").Value = $newTimestamp
$ws.Range("C8").Value = $newTimestamp
$ws.Range("C9").Value = $newTimestamp
$ws.Range("C10").Value = $newTimestamp
$ws.Range("C11").Value = $newTimestamp
$ws.Range("C12").Value = $newTimestamp
$ws.Range("C13").Value = $newTimestamp
$ws.Range("C14").Value = $newTimestamp
$ws.Range("C15").Value = $newTimestamp
